$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Cells receiving plain numeric-looking text need to be forced to Text format
# first, otherwise Excel will silently convert them to numbers and the exact
# string representation (trailing zeros, etc.) would be lost.
$textCells = "D5,D6,D8,D9,D10,D11,D13,D14,D15,D16,D20,D21,D22,D23,D24,D25,D26,D27,D28,D29,D30,D31,D32,D33,D34,D35,D36,D37,D38,D39,D40,D41,D43,D44,D45,D47,D48,D49,D50,D51" -split ","
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "29.945.15"
$ws.Range("E2").Value = "  +0.45%  "
$ws.Range("D3").Value = "1.892.94"
$ws.Range("E3").Value = "  +0.02%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "0.7791"
$ws.Range("E5").Value = "  +0.30%  "
$ws.Range("D6").Value = "244.00"
$ws.Range("E6").Value = "  +0.06%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").Value = "0.3137"
$ws.Range("E8").Value = "  +0.41%  "
$ws.Range("D9").Value = "25.86"
$ws.Range("E9").Value = "  +2.33%  "
$ws.Range("D10").Value = "0.07280"
$ws.Range("E10").Value = "  +1.51%  "
$ws.Range("B11").Value = "TRON"
$ws.Range("C11").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D11").Value = "0.08694"
$ws.Range("E11").Value = "  +7.85%  "
$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").Value = "2.031.51"
$ws.Range("E12").Value = "  +8.20%  "
$ws.Range("D13").Value = "0.7750"
$ws.Range("E13").Value = "  +1.27%  "
$ws.Range("D14").Value = "5.414"
$ws.Range("E14").Value = "  -0.65%  "
$ws.Range("D15").Value = "94.55"
$ws.Range("E15").Value = "  +2.50%  "
$ws.Range("D16").Value = "6.207"
$ws.Range("E16").Value = "  +0.67%  "
$ws.Range("D17").Value = "29.950.19"
$ws.Range("E17").Value = "  +0.53%  "
$ws.Range("E18").Value = "  -0.11%  "
$ws.Range("B19").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C19").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D19").Value = "2.319.62"
$ws.Range("E19").Value = "  +9.89%  "
$ws.Range("B20").Value = "BitcoinCash"
$ws.Range("C20").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D20").Value = "246.08"
$ws.Range("E20").Value = "  +1.14%  "
$ws.Range("D21").Value = "0.000007884"
$ws.Range("E21").Value = "  +1.61%  "
$ws.Range("D22").Value = "8.275"
$ws.Range("E22").Value = "  +2.09%  "
$ws.Range("D23").Value = "1.000"
$ws.Range("E23").Value = "  +0.08%  "
$ws.Range("D24").Value = "1.000"
$ws.Range("D25").Value = "0.1673"
$ws.Range("E25").Value = "  +4.39%  "
$ws.Range("D26").Value = "9.508"
$ws.Range("E26").Value = "  +1.19%  "
$ws.Range("D27").Value = "163.36"
$ws.Range("E27").Value = "  +1.15%  "
$ws.Range("D28").Value = "18.87"
$ws.Range("E28").Value = "  +0.75%  "
$ws.Range("D29").Value = "2.052"
$ws.Range("E29").Value = "  +0.21%  "
$ws.Range("D30").Value = "1.432"
$ws.Range("E30").Value = "  +0.08%  "
$ws.Range("D31").Value = "1.542"
$ws.Range("E31").Value = "  -0.36%  "
$ws.Range("D32").Value = "4.516"
$ws.Range("E32").Value = "  +1.07%  "
$ws.Range("D33").Value = "4.132"
$ws.Range("E33").Value = "  +0.85%  "
$ws.Range("D34").Value = "0.05491"
$ws.Range("E34").Value = "  -0.69%  "
$ws.Range("D35").Value = "1.248"
$ws.Range("E35").Value = "  -1.20%  "
$ws.Range("D36").Value = "0.7569"
$ws.Range("E36").Value = "  +1.61%  "
$ws.Range("D37").Value = "1.006"
$ws.Range("E37").Value = "  +0.86%  "
$ws.Range("D38").Value = "2.686"
$ws.Range("E38").Value = "  +2.61%  "
$ws.Range("D39").Value = "0.01960"
$ws.Range("E39").Value = "  +2.51%  "
$ws.Range("D40").Value = "2.792"
$ws.Range("E40").Value = "  +0.52%  "
$ws.Range("D41").Value = "0.4518"
$ws.Range("E41").Value = "  +2.23%  "
$ws.Range("D42").Value = "1.113.14"
$ws.Range("E42").Value = "  -2.29%  "
$ws.Range("D43").Value = "74.15"
$ws.Range("E43").Value = "  +0.77%  "
$ws.Range("D44").Value = "6.092"
$ws.Range("E44").Value = "  +4.10%  "
$ws.Range("D45").Value = "0.8532"
$ws.Range("E45").Value = "  +0.01%  "
$ws.Range("B46").Value = "RocketPoolETH"
$ws.Range("C46").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D46").Value = "2.205.14"
$ws.Range("E46").Value = "  +9.81%  "
$ws.Range("B47").Value = "PaxDollar"
$ws.Range("C47").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D47").Value = "0.9998"
$ws.Range("E47").Value = "  -0.03%  "
$ws.Range("B48").Value = "Quant"
$ws.Range("C48").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D48").Value = "103.41"
$ws.Range("E48").Value = "  -0.38%  "
$ws.Range("B49").Value = "RenderToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D49").Value = "1.887"
$ws.Range("E49").Value = "  -0.01%  "
$ws.Range("D50").Value = "7.607"
$ws.Range("E50").Value = "  +2.26%  "
$ws.Range("D51").Value = "9.888"
$ws.Range("E51").Value = "  -0.21%  "

# Restore the default (unstyled) appearance now that the text values are locked in
foreach ($addr in $textCells) {
    $ws.Range($addr).Style = "Normal"
}
